$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("D2").Value = "'67.746.19"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.24%  "
$ws.Range("D3").Value = "'3.326.89"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +1.01%  "
$ws.Range("D4").Value = "'1.00"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").Value = "'582.45"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.30%  "
$ws.Range("D6").Value = "'174.58"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.08%  "
$ws.Range("D7").Value = "'1.00"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.04%  "
$ws.Range("E8").Value = "  +1.58%  "
$ws.Range("D9").Value = "'3.320.62"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +1.04%  "
$ws.Range("D10").Value = "'0.181"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +4.62%  "
$ws.Range("E11").Value = "  +1.34%  "
$ws.Range("D12").Value = "'46.97"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +4.09%  "
$ws.Range("D13").Value = "'0.0000271"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +1.70%  "
$ws.Range("D14").Value = "'697.72"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +4.90%  "
$ws.Range("D15").Value = "'3.867.53"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.00%  "
$ws.Range("E16").Value = "  +0.69%  "
$ws.Range("D17").Value = "'67.790.14"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.14%  "
$ws.Range("E18").Value = "  +0.62%  "
$ws.Range("D19").Value = "'3.302.95"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.05%  "
$ws.Range("E20").Value = "  +0.63%  "
$ws.Range("E21").Value = "  +2.69%  "
$ws.Range("E22").Value = "  +0.76%  "
$ws.Range("E23").Value = "  +0.34%  "
$ws.Range("D24").Value = "'16.90"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.32%  "
$ws.Range("D25").Value = "'101.27"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +3.12%  "
$ws.Range("D26").Value = "'3.90"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +1.46%  "
$ws.Range("E27").Value = "  +1.30%  "
$ws.Range("E28").Value = "  +2.80%  "
$ws.Range("D29").Value = "'32.84"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.37%  "
$ws.Range("D30").Value = "'8.51"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +2.30%  "
$ws.Range("D31").Value = "'6.97"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.31%  "
$ws.Range("D32").Value = "'571.74"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.52%  "
$ws.Range("E33").Value = "  +0.82%  "
$ws.Range("E34").Value = "  +2.25%  "
$ws.Range("B35").Value = "Maker"
$ws.Range("C35").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D35").Value = "'3.718.87"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.96%  "
$ws.Range("B36").Value = "Dai"
$ws.Range("C36").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D36").Value = "'1.00"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.19%  "
$ws.Range("D37").Value = "'56.47"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +1.55%  "
$ws.Range("D38").Value = "'3.27"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -2.55%  "
$ws.Range("D39").Value = "'35.61"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +10.66%  "
$ws.Range("E40").Value = "  +2.85%  "
$ws.Range("D41").Value = "'3.14"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +2.83%  "
$ws.Range("E42").Value = "  -0.29%  "
$ws.Range("B43").Value = "PEPE"
$ws.Range("C43").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D43").Value = "'0.0₃0669"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +1.48%  "
$ws.Range("B44").Value = "TheGraph"
$ws.Range("C44").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D44").Value = "'0.334"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +2.39%  "
$ws.Range("B45").Value = "ApeXProtocol"
$ws.Range("C45").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D45").Value = "'3.30"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +2.39%  "
$ws.Range("E46").Value = "  +0.91%  "
$ws.Range("D47").Value = "'2.62"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +1.80%  "
$ws.Range("E48").Value = "  +1.44%  "
$ws.Range("E49").Value = "  -0.08%  "
$ws.Range("E50").Value = "  -2.37%  "
$ws.Range("D51").Value = "'131.23"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +1.88%  "
